$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.102.34'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.822.53'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4626'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.46%  '
$ws.Range('E8').Value = '  -1.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07305'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8697'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.13'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.873.57'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07580'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.348'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.61'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.476'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.86%  '
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008646'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.35%  '
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.351.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('E22').Value = '  -2.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.095.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.876'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.095'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.07%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.095'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08911'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.959'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7337'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.74%  '
$ws.Range('E34').Value = '  -2.58%  '
$ws.Range('E35').Value = '  -3.35%  '
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.484'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05258'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.069'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01920'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.929'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.160'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5210'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1630'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.265'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4885'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.010'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '103.75'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.632'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06254'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.31%  '
